$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 11
$ws.Range("D2").Value = 0.606
$ws.Range("E2").Value = 0.805
$ws.Range("F2").Value = 0.798
$ws.Range("G2").Value = 0.232
$ws.Range("H2").Value = 1.528

$ws.Range("C3").Value = 6
$ws.Range("D3").Value = 0.606
$ws.Range("E3").Value = 0.966
$ws.Range("F3").Value = 0.893
$ws.Range("G3").Value = 0.22
$ws.Range("H3").Value = 1.708

$ws.Range("C4").Value = 8
$ws.Range("D4").Value = 0.613
$ws.Range("E4").Value = 0.984
$ws.Range("F4").Value = 1.037
$ws.Range("G4").Value = 0.235
$ws.Range("H4").Value = 1.609

$ws.Range("D5").Value = 0.762
$ws.Range("E5").Value = 0.75
$ws.Range("F5").Value = 0.75
$ws.Range("G5").Value = 0.47
$ws.Range("H5").Value = 2.074

$ws.Range("D6").Value = 0.618
$ws.Range("E6").Value = 0.403
$ws.Range("F6").Value = 0.419
$ws.Range("G6").Value = 0.656
$ws.Range("H6").Value = 1.533

$ws.Range("C7").Value = 3
$ws.Range("D7").Value = 1.021
$ws.Range("E7").Value = 1.405
$ws.Range("F7").Value = 1.491
$ws.Range("G7").Value = 0.019
$ws.Range("H7").Value = 2.356
